$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 79882.836
$ws.Range("J3").Value = 79882.836
$ws.Range("L3").Value = 79882.836
$ws.Range("N3").Value = -80110.836
$ws.Range("H12").Value = 255.81818
$ws.Range("I12").Value = 256
$ws.Range("K12").Value = 256
$ws.Range("M12").Value = -86
$ws.Range("H39").Value = 430.31033
$ws.Range("I39").Value = 71.375
$ws.Range("J39").Value = 567.0476
$ws.Range("K39").Value = 214.125
$ws.Range("L39").Value = 1701.1428
$ws.Range("M39").Value = 81.875
$ws.Range("N39").Value = -2293.1428
$ws.Range("H40").Value = 3863.5
$ws.Range("I40").Value = 2424.5
$ws.Range("K40").Value = 2424.5
$ws.Range("M40").Value = -2249.5
$ws.Range("H42").Value = 226.57143
$ws.Range("I42").Value = 226.57143
$ws.Range("K42").Value = 679.71429
$ws.Range("M42").Value = -449.71429
$ws.Range("H57").Value = 30319.75
$ws.Range("J57").Value = 30319.75
$ws.Range("L57").Value = 90959.25
$ws.Range("N57").Value = -91957.25
$ws.Range("H92").Value = 528.7692
$ws.Range("I92").Value = 505.38095
$ws.Range("K92").Value = 505.38095
$ws.Range("M92").Value = 742.61905
$ws.Range("H100").Value = 5439.9
$ws.Range("I100").Value = 5928
$ws.Range("K100").Value = 5928
$ws.Range("M100").Value = -5387
$ws.Range("H102").Value = 79882.836
$ws.Range("J102").Value = 79882.836
$ws.Range("L102").Value = 79882.836
$ws.Range("N102").Value = -86372.836
$ws.Range("H111").Value = 2640.875
$ws.Range("I111").Value = 2640.875
$ws.Range("K111").Value = 7922.625
$ws.Range("M111").Value = -4855.625
$ws.Range("H116").Value = 6867.9375
$ws.Range("J116").Value = 6876.3335
$ws.Range("L116").Value = 6876.3335
$ws.Range("N116").Value = -13760.3335
$ws.Range("H128").Value = 88000
$ws.Range("J128").Value = 88000
$ws.Range("L128").Value = 88000
$ws.Range("N128").Value = -97960
$ws.Range("H137").Value = 9558.166999999999
$ws.Range("I137").Value = 2650
$ws.Range("J137").Value = 13012.25
$ws.Range("K137").Value = 7950
$ws.Range("L137").Value = 39036.75
$ws.Range("M137").Value = -5400
$ws.Range("N137").Value = -44136.75
$ws.Range("H138").Value = 2724.0857
$ws.Range("I138").Value = 940.1111
$ws.Range("J138").Value = 3341.6155
$ws.Range("K138").Value = 2820.3333
$ws.Range("L138").Value = 10024.8465
$ws.Range("M138").Value = 2319.6667
$ws.Range("N138").Value = -20304.8465

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2029.5
$ws.Range("I45").Value = 2126.3333
$ws.Range("K45").Value = 2126.3333
$ws.Range("M45").Value = -1749.3333
$ws.Range("H76").Value = 80000
$ws.Range("J76").Value = 80000
$ws.Range("L76").Value = 80000
$ws.Range("N76").Value = -80676
$ws.Range("H79").Value = 80000
$ws.Range("J79").Value = 80000
$ws.Range("L79").Value = 80000
$ws.Range("N79").Value = -82340
$ws.Range("H102").Value = 8862.1875
$ws.Range("I102").Value = 8715.842000000001
$ws.Range("K102").Value = 8715.842000000001
$ws.Range("M102").Value = -7093.842000000001
$ws.Range("H119").Value = 75000.5
$ws.Range("J119").Value = 75000.5
$ws.Range("L119").Value = 75000.5
$ws.Range("N119").Value = -84676.5
$ws.Range("H122").Value = 1242.3334
$ws.Range("I122").Value = 1085.125
$ws.Range("K122").Value = 3255.375
$ws.Range("M122").Value = -805.375
$ws.Range("H132").Value = 6006.029
$ws.Range("I132").Value = 3436.7693
$ws.Range("J132").Value = 13428.333
$ws.Range("K132").Value = 10310.3079
$ws.Range("L132").Value = 40284.999
$ws.Range("M132").Value = -7780.3079
$ws.Range("N132").Value = -45344.999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 30224.666
$ws.Range("J96").Value = 65208
$ws.Range("L96").Value = 65208
$ws.Range("N96").Value = -70700

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 55285.715
$ws.Range("I59").Value = 12000
$ws.Range("K59").Value = 12000
$ws.Range("M59").Value = -10855
$ws.Range("H62").Value = 4999.5
$ws.Range("J62").Value = 5999
$ws.Range("L62").Value = 5999
$ws.Range("N62").Value = -7247
$ws.Range("H65").Value = 4999.5
$ws.Range("J65").Value = 5999
$ws.Range("L65").Value = 29995
$ws.Range("N65").Value = -36235
$ws.Range("H132").Value = 2882.0322
$ws.Range("I132").Value = 2379.9
$ws.Range("J132").Value = 3795
$ws.Range("K132").Value = 7139.700000000001
$ws.Range("L132").Value = 11385
$ws.Range("M132").Value = -4609.700000000001
$ws.Range("N132").Value = -16445

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5999.8
$ws.Range("I80").Value = 5999
$ws.Range("K80").Value = 17997
$ws.Range("M80").Value = -17061
$ws.Range("H83").Value = 5999.8
$ws.Range("I83").Value = 5999
$ws.Range("K83").Value = 53991
$ws.Range("M83").Value = -49311
$ws.Range("H107").Value = 484
$ws.Range("I107").Value = 285
$ws.Range("J107").Value = 508.875
$ws.Range("K107").Value = 855
$ws.Range("L107").Value = 1526.625
$ws.Range("M107").Value = 1065
$ws.Range("N107").Value = -5366.625
$ws.Range("H137").Value = 4921.143
$ws.Range("J137").Value = 3466.8572
$ws.Range("L137").Value = 10400.5716
$ws.Range("N137").Value = -20600.5716

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4712.2856
$ws.Range("I70").Value = 4712.2856
$ws.Range("K70").Value = 4712.2856
$ws.Range("M70").Value = -4442.2856
$ws.Range("H73").Value = 4712.2856
$ws.Range("I73").Value = 4712.2856
$ws.Range("K73").Value = 4712.2856
$ws.Range("M73").Value = -3776.2856
$ws.Range("H136").Value = 26849.875
$ws.Range("J136").Value = 26849.875
$ws.Range("L136").Value = 80549.625
$ws.Range("N136").Value = -85649.625

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49661.137
$ws.Range("I7").Value = 3472.1177
$ws.Range("J7").Value = 206703.8
$ws.Range("K7").Value = 3472.1177
$ws.Range("L7").Value = 206703.8
$ws.Range("M7").Value = -3360.1177
$ws.Range("N7").Value = -206927.8
$ws.Range("H40").Value = 2422.1482
$ws.Range("I40").Value = 1799.9565
$ws.Range("J40").Value = 5999.75
$ws.Range("K40").Value = 1799.9565
$ws.Range("L40").Value = 5999.75
$ws.Range("M40").Value = -1663.9565
$ws.Range("N40").Value = -6271.75
$ws.Range("H126").Value = 49661.137
$ws.Range("I126").Value = 3472.1177
$ws.Range("J126").Value = 206703.8
$ws.Range("K126").Value = 10416.3531
$ws.Range("L126").Value = 620111.3999999999
$ws.Range("M126").Value = -7946.3531
$ws.Range("N126").Value = -625051.3999999999
$ws.Range("H132").Value = 480043.7
$ws.Range("I132").Value = 456409.1
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 1369227.3
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -1366697.3
$ws.Range("N132").Value = -3005075
$ws.Range("H136").Value = 104828.91
$ws.Range("I136").Value = 5311.3
$ws.Range("J136").Value = 1100005
$ws.Range("K136").Value = 15933.9
$ws.Range("L136").Value = 3300015
$ws.Range("M136").Value = -13383.9
$ws.Range("N136").Value = -3305115

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 543339.7
$ws.Range("I31").Value = 1500000
$ws.Range("J31").Value = 65009.5
$ws.Range("K31").Value = 1500000
$ws.Range("L31").Value = 65009.5
$ws.Range("M31").Value = -1499652
$ws.Range("N31").Value = -65705.5
$ws.Range("H81").Value = 1489.5
$ws.Range("I81").Value = 1489.5
$ws.Range("K81").Value = 2979
$ws.Range("M81").Value = -1918
$ws.Range("H84").Value = 1489.5
$ws.Range("I84").Value = 1489.5
$ws.Range("K84").Value = 14895
$ws.Range("M84").Value = -9591
$ws.Range("H119").Value = 85198
$ws.Range("J119").Value = 85198
$ws.Range("L119").Value = 85198
$ws.Range("N119").Value = -94874
$ws.Range("H132").Value = 1951.85
$ws.Range("I132").Value = 1951.85
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5855.549999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3325.549999999999
$ws.Range("N132").ClearContents()
